# Quarterly indexing esoteric bug-fix operation
#
# Column A holds (mis-indexed) quarter-start date serials. The fix shifts
# each date one month forward and re-anchors it to the 15th of that month
# (i.e. old = 1st-of-quarter -> new = 15th of the following month), with a
# proper December -> January / year+1 rollover.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 73
$col = 1   # column A

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $old = $cell.Value2

    if ($old -ne $null) {
        $d = [DateTime]::FromOADate($old)

        $newMonth = $d.Month + 1
        $newYear = $d.Year
        if ($newMonth -gt 12) {
            $newMonth = $newMonth - 12
            $newYear = $newYear + 1
        }

        $d2 = Get-Date -Year $newYear -Month $newMonth -Day 15 -Hour 0 -Minute 0 -Second 0 -Millisecond 0
        $new = [Math]::Floor($d2.ToOADate())

        $cell.Value = $new
    }
}
